$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price/Volume columns so Excel does not
# auto-convert numeric-looking strings (e.g. "1.525", "134.00") into
# numbers, which would lose formatting / precision vs. the source text.
$priceVolRange = $ws.Range("D2:E51")
$priceVolRange.NumberFormat = "@"

$ws.Range("D2").Value = "24.551.29"
$ws.Range("D3").Value = "1.693.43"
$ws.Range("E3").Value = "  +1.60%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "314.93"
$ws.Range("E5").Value = "  +1.61%  "
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("D7").Value = "0.3932"
$ws.Range("E7").Value = "  +1.37%  "
$ws.Range("D8").Value = "0.3995"
$ws.Range("E8").Value = "  +1.10%  "
$ws.Range("D9").Value = "1.525"
$ws.Range("E9").Value = "  +4.95%  "
$ws.Range("E10").Value = "  +0.15%  "
$ws.Range("D11").Value = "51.96"
$ws.Range("E11").Value = "  +2.34%  "
$ws.Range("D12").Value = "0.08728"
$ws.Range("E12").Value = "  +0.87%  "
$ws.Range("D13").Value = "7.205"
$ws.Range("E13").Value = "  +6.58%  "
$ws.Range("D14").Value = "23.08"
$ws.Range("E14").Value = "  +2.11%  "
$ws.Range("D15").Value = "0.00001314"
$ws.Range("E15").Value = "  -0.06%  "
$ws.Range("D16").Value = "7.578"
$ws.Range("E16").Value = "  +4.06%  "
$ws.Range("D17").Value = "1.691.08"
$ws.Range("E17").Value = "  +1.87%  "
$ws.Range("D18").Value = "99.55"
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("D19").Value = "0.07057"
$ws.Range("E19").Value = "  +3.93%  "
$ws.Range("D20").Value = "19.57"
$ws.Range("E20").Value = "  +2.35%  "
$ws.Range("D21").Value = "6.865"
$ws.Range("E21").Value = "  +3.22%  "
$ws.Range("E22").Value = "  +0.31%  "
$ws.Range("D23").Value = "14.03"
$ws.Range("D24").Value = "24.552.45"
$ws.Range("E24").Value = "  +3.19%  "
$ws.Range("D25").Value = "3.076"
$ws.Range("E25").Value = "  +8.12%  "
$ws.Range("D26").Value = "2.327"
$ws.Range("E26").Value = "  +0.55%  "
$ws.Range("D27").Value = "22.29"
$ws.Range("E27").Value = "  +2.50%  "
$ws.Range("D28").Value = "160.86"
$ws.Range("E28").Value = "  +0.95%  "
$ws.Range("D29").Value = "5.215"
$ws.Range("E29").Value = "  +0.79%  "
$ws.Range("D30").Value = "134.00"
$ws.Range("E30").Value = "  +3.31%  "
$ws.Range("D31").Value = "7.544"
$ws.Range("E31").Value = "  +11.92%  "
$ws.Range("D32").Value = "1.876.74"
$ws.Range("E32").Value = "  +1.21%  "
$ws.Range("D33").Value = "1.088"
$ws.Range("E33").Value = "  -2.78%  "
$ws.Range("D34").Value = "0.08546"
$ws.Range("E34").Value = "  +0.78%  "
$ws.Range("D35").Value = "7.272"
$ws.Range("E35").Value = "  +9.13%  "
$ws.Range("D36").Value = "11.29"
$ws.Range("E36").Value = "  +8.07%  "
$ws.Range("D37").Value = "1.943"
$ws.Range("E37").Value = "  -1.46%  "
$ws.Range("D38").Value = "0.2704"
$ws.Range("E38").Value = "  +1.79%  "
$ws.Range("D39").Value = "14.41"
$ws.Range("E39").Value = "  -0.57%  "
$ws.Range("D40").Value = "0.02738"
$ws.Range("E40").Value = "  +8.83%  "
$ws.Range("D41").Value = "0.09009"
$ws.Range("E41").Value = "  +2.52%  "
$ws.Range("D42").Value = "1.468"
$ws.Range("E42").Value = "  +0.62%  "
$ws.Range("D43").Value = "0.7656"
$ws.Range("E43").Value = "  +1.02%  "
$ws.Range("D44").Value = "0.7157"
$ws.Range("E44").Value = "  +1.72%  "
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").Value = "2.522"
$ws.Range("E45").Value = "  +3.96%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "15.27"
$ws.Range("E46").Value = "  +1.46%  "
$ws.Range("D47").Value = "4.199"
$ws.Range("E47").Value = "  +2.31%  "
$ws.Range("D48").Value = "1.000"
$ws.Range("E48").Value = "  +0.15%  "
$ws.Range("B49").Value = "Flow"
$ws.Range("C49").Value = "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"
$ws.Range("D49").Value = "1.330"
$ws.Range("E49").Value = "  +8.61%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "140.64"
$ws.Range("E50").Value = "  +1.11%  "
$ws.Range("D51").Value = "0.07992"

# Restore the default cell style so we do not leave a stray number
# format applied (matches original workbook styling).
$priceVolRange.Style = "Normal"
